# Adds two new players ("Anne-Lise" and "Jean Rob") to the list,
# giving Didier the possibility to specify his killer.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("A10").Value = "Anne-Lise"
$ws.Range("A11").Value = "Jean Rob"

$ws.Range("A11").Select()
